$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AV2: ID -> RLN
$ws.Range("AV2").Value = "RLN"

# AR7: - -> ZL
$ws.Range("AR7").Value = "ZL"

# Row 13: VMW/TC -> VMW/TAC
$ws.Range("T13").Value = "VMW/TAC"
$ws.Range("X13").Value = "VMW/TAC"
$ws.Range("AB13").Value = "VMW/TAC"
$ws.Range("AF13").Value = "VMW/TAC"
$ws.Range("AJ13").Value = "VMW/TAC"
$ws.Range("AN13").Value = "VMW/TAC"

# Row 14: - -> AMS/CS/IMB/JML
$ws.Range("T14").Value = "AMS/CS/IMB/JML"
$ws.Range("X14").Value = "AMS/CS/IMB/JML"
$ws.Range("AB14").Value = "AMS/CS/IMB/JML"
$ws.Range("AF14").Value = "AMS/CS/IMB/JML"
$ws.Range("AJ14").Value = "AMS/CS/IMB/JML"
$ws.Range("AN14").Value = "AMS/CS/IMB/JML"

# Row 16: IK/MS -> MS/MC ; MC/SD -> SD/MC
$ws.Range("T16").Value = "MS/MC"
$ws.Range("X16").Value = "MS/MC"
$ws.Range("AB16").Value = "MS/MC"
$ws.Range("AF16").Value = "SD/MC"
$ws.Range("AJ16").Value = "SD/MC"
$ws.Range("AN16").Value = "SD/MC"

# Row 19: JW -> JW/TM or JW/DK
$ws.Range("T19").Value = "JW/TM"
$ws.Range("X19").Value = "JW/DK"
$ws.Range("AB19").Value = "JW/TM"
$ws.Range("AF19").Value = "JW/DK"
$ws.Range("AJ19").Value = "JW/TM"
$ws.Range("AN19").Value = "JW/TM"

# Row 20: SC,DK -> SC/DK ; SS,MK -> SS/MK
$ws.Range("T20").Value = "SC/DK"
$ws.Range("X20").Value = "SC/DK"
$ws.Range("AB20").Value = "SC/DK"
$ws.Range("AF20").Value = "SS/MK"
$ws.Range("AJ20").Value = "SS/MK"
$ws.Range("AN20").Value = "SS/MK"

# Row 21: LAB 1 -> numeric 62 / 61
$ws.Range("U21").Value = 62
$ws.Range("Y21").Value = 62
$ws.Range("AC21").Value = 62
$ws.Range("AG21").Value = 61
$ws.Range("AK21").Value = 61
$ws.Range("AO21").Value = 61

# Row 22: FAL,BM,JNS -> FAL/BM/JNS ; ADM,DL,KM -> ADM/DL/KM
$ws.Range("T22").Value = "FAL/BM/JNS"
$ws.Range("X22").Value = "FAL/BM/JNS"
$ws.Range("AB22").Value = "FAL/BM/JNS"
$ws.Range("AF22").Value = "ADM/DL/KM"
$ws.Range("AJ22").Value = "ADM/DL/KM"
$ws.Range("AN22").Value = "ADM/DL/KM"

# Row 25: Fr. GN / FR GN / FR GN/CK / GN/CK -> FGN / FGN/CK
$ws.Range("T25").Value = "FGN"
$ws.Range("X25").Value = "FGN"
$ws.Range("AB25").Value = "FGN/CK"
$ws.Range("AF25").Value = "FGN/CK"
$ws.Range("AJ25").Value = "FGN/CK"
$ws.Range("AN25").Value = "FGN/CK"
